# "final fixes before submiting"
#
# Habit2's sample data (x/y points feeding the scatter chart) was
# refreshed with a new run, the series got renamed from the generic
# "Results" to "ApplicationCount", and Habit2 was left as the active /
# selected sheet (with the cursor sitting on L12) instead of Habit3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Habit2")

# --- series label -----------------------------------------------------
$ws.Range("A2").Value = "ApplicationCount"

# --- x values (row 1) ---------------------------------------------------
$ws.Range("B1").Value = 0.01
$ws.Range("C1").Value = 0.025
$ws.Range("D1").Value = 0.05
$ws.Range("E1").Value = 0.075
# F1:K1 (0.1, 0.2, 0.25, 0.3, 0.35, 0.4) are unchanged.

# --- y values (row 2) ----------------------------------------------------
$ws.Range("B2").Value = 0.956
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.988
$ws.Range("F2").Value = 0.08
# G2:K2 (all 1) are unchanged.

# Keep the chart's cached series text/number caches in step with the
# refreshed cells above.
$co = $ws.ChartObjects(1)
$series = $co.Chart.SeriesCollection(1)
$series.XValues = $ws.Range("B1:K1")
$series.Values = $ws.Range("B2:K2")

# --- window/selection state --------------------------------------------
# Habit2 becomes the active tab, with L12 selected (this also clears
# tabSelected on the previously-active Habit3 automatically).
$ws.Range("L12").Select()
